$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update serial number (kept as text), customer name, and clear the
# broken/duplicate URLs that were causing viewer page errors.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "56888"
$ws.Range("D2").Value = "Bay 1 flightline 2"
$ws.Range("F2").ClearContents()

# Row 3: update flightline, serial number, customer name, and clear the
# test/broken URLs.
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = "Test"
$ws.Range("D3").Value = "Test"
$ws.Range("F3").ClearContents()
